# Wrapping up edits to Ch4 (DEBkiss results / Figure 1 chapter 4)
#
# 1. Refresh the cached "datetimeFigureOut" field text (slide master + all
#    11 custom layouts) from 3/22/2023 -> 4/2/2023.
# 2. Nudge four existing shapes on slide 1 (two labels + their red outline
#    rectangles) to their new positions.
# 3. Add two new bold "section title" textboxes to slide 1.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholders
# ---------------------------------------------------------------------------
$newDate = "4/2/2023"

function Update-DateShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "*Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DateShapes $p.SlideMaster
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateShapes $p.SlideMaster.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------------
# 2. Reposition the four shapes on slide 1
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# TextBox 22 ("Embryo mortality rate (...)") -> moved left/down
$shp23 = Get-ShapeById $s 23
$shp23.Left = 750.7703247070312
$shp23.Top = 134.33291625976562

# TextBox 23 ("Post-hatch mortality rate (...)") -> moved right
$shp24 = Get-ShapeById $s 24
$shp24.Left = 895.5204467773438
$shp24.Top = 393.31805419921875

# Rectangle 31 (red outline around TextBox 23) -> moved right, same as above
$shp32 = Get-ShapeById $s 32
$shp32.Left = 895.5203247070312
$shp32.Top = 392.3668518066406

# Rectangle 32 (red outline around TextBox 22) -> moved left/down, same as above
$shp33 = Get-ShapeById $s 33
$shp33.Left = 750.0455322265625
$shp33.Top = 134.0693817138672

# ---------------------------------------------------------------------------
# 3. Add the two new section-title textboxes
# ---------------------------------------------------------------------------
$tb6 = $s.Shapes.AddTextbox(1, 204.2975616455078, 73.5394515991211, 304.0155944824219, 31.504724502563477)
$tb6.TextFrame.WordWrap = 1
$tb6.TextFrame.AutoSize = 1
$tr6 = $tb6.TextFrame.TextRange
$tr6.Text = "Full Life Cycle Energy Budget"
$tr6.Font.Name = "Arial"
$tr6.Font.NameComplexScript = "Arial"
$tr6.Font.Size = 20
$tr6.Font.Bold = 1
$tb6.Fill.Visible = 0

$tb7 = $s.Shapes.AddTextbox(1, 806.8411865234375, 73.1900863647461, 244.82504272460938, 31.504724502563477)
$tb7.TextFrame.WordWrap = 1
$tb7.TextFrame.AutoSize = 1
$tr7 = $tb7.TextFrame.TextRange
$tr7.Text = "Stage Specific Survival"
$tr7.Font.Name = "Arial"
$tr7.Font.NameComplexScript = "Arial"
$tr7.Font.Size = 20
$tr7.Font.Bold = 1
$tb7.Fill.Visible = 0

Write-Output "done"
